# Scheduled-runner market-data refresh: updates the derived price/profit
# columns (H:N) on 36 leve rows across all 8 job sheets to match the
# latest Universalis price snapshot. Generated from the authoritative diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1942.875
$ws.Range("I100").Value = 1671.4546
$ws.Range("J100").Value = 2540
$ws.Range("K100").Value = 1671.4546
$ws.Range("L100").Value = 2540
$ws.Range("M100").Value = -1130.4546
$ws.Range("N100").Value = -3622

$ws.Range("H134").Value = 47236.75
$ws.Range("J134").Value = 47236.75
$ws.Range("L134").Value = 47236.75
$ws.Range("N134").Value = -57376.75

$ws.Range("H137").Value = 1590.4445
$ws.Range("I137").Value = 1646.9474
$ws.Range("J137").Value = 1456.25
$ws.Range("K137").Value = 4940.8422
$ws.Range("L137").Value = 4368.75
$ws.Range("M137").Value = -2390.8422
$ws.Range("N137").Value = -9468.75

$ws.Range("H138").Value = 2445.2368
$ws.Range("I138").Value = 1396.3158
$ws.Range("J138").Value = 3494.158
$ws.Range("K138").Value = 4188.9474
$ws.Range("L138").Value = 10482.474
$ws.Range("M138").Value = 951.0526
$ws.Range("N138").Value = -20762.474

$ws.Range("I141").Value = 2926.2856
$ws.Range("J141").Value = 3267.5
$ws.Range("K141").Value = 8778.856800000001
$ws.Range("L141").Value = 9802.5
$ws.Range("M141").Value = -3598.856800000001
$ws.Range("N141").Value = -20162.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 635.92
$ws.Range("I2").Value = 690.0952
$ws.Range("J2").Value = 351.5
$ws.Range("K2").Value = 690.0952
$ws.Range("L2").Value = 351.5
$ws.Range("M2").Value = -577.0952
$ws.Range("N2").Value = -577.5

$ws.Range("H32").Value = 5592.459
$ws.Range("I32").Value = 4350.5796
$ws.Range("J32").Value = 10948.0625
$ws.Range("K32").Value = 4350.5796
$ws.Range("L32").Value = 10948.0625
$ws.Range("M32").Value = -4063.5796
$ws.Range("N32").Value = -11522.0625

$ws.Range("H74").Value = 22729548
$ws.Range("I74").Value = 33335252
$ws.Range("J74").Value = 3036.5715
$ws.Range("K74").Value = 33335252
$ws.Range("L74").Value = 3036.5715
$ws.Range("M74").Value = -33334378
$ws.Range("N74").Value = -4784.5715

$ws.Range("H77").Value = 22729548
$ws.Range("I77").Value = 33335252
$ws.Range("J77").Value = 3036.5715
$ws.Range("K77").Value = 166676260
$ws.Range("L77").Value = 15182.8575
$ws.Range("M77").Value = -166671892
$ws.Range("N77").Value = -23918.8575

$ws.Range("H102").Value = 989.36365
$ws.Range("I102").Value = 860.7895
$ws.Range("K102").Value = 860.7895
$ws.Range("M102").Value = 761.2105

$ws.Range("H112").Value = 29515.834
$ws.Range("J112").Value = 29515.834
$ws.Range("L112").Value = 29515.834
$ws.Range("N112").Value = -32469.834

$ws.Range("H116").Value = 635.92
$ws.Range("I116").Value = 690.0952
$ws.Range("J116").Value = 351.5
$ws.Range("K116").Value = 690.0952
$ws.Range("L116").Value = 351.5
$ws.Range("M116").Value = 1603.9048
$ws.Range("N116").Value = -4939.5

$ws.Range("H122").Value = 1740.1333
$ws.Range("I122").Value = 1793
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 5379
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -2929
$ws.Range("N122").Value = -7900

$ws.Range("H125").Value = 31994.5
$ws.Range("J125").Value = 31994.5
$ws.Range("L125").Value = 31994.5
$ws.Range("N125").Value = -41834.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 635.92
$ws.Range("I3").Value = 690.0952
$ws.Range("J3").Value = 351.5
$ws.Range("K3").Value = 690.0952
$ws.Range("L3").Value = 351.5
$ws.Range("M3").Value = -576.0952
$ws.Range("N3").Value = -579.5

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H20").Value = 3328.4546
$ws.Range("I20").Value = 3872
$ws.Range("J20").Value = 2377.25
$ws.Range("K20").Value = 3872
$ws.Range("L20").Value = 2377.25
$ws.Range("M20").Value = -3625
$ws.Range("N20").Value = -2871.25

$ws.Range("H107").Value = 1189.8334
$ws.Range("I107").Value = 957.36365
$ws.Range("K107").Value = 957.36365
$ws.Range("M107").Value = 962.63635

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3918.0334
$ws.Range("I31").Value = 968.1818
$ws.Range("J31").Value = 5625.8423
$ws.Range("K31").Value = 968.1818
$ws.Range("L31").Value = 5625.8423
$ws.Range("M31").Value = -673.1818
$ws.Range("N31").Value = -6215.8423

$ws.Range("H34").Value = 3918.0334
$ws.Range("I34").Value = 968.1818
$ws.Range("J34").Value = 5625.8423
$ws.Range("K34").Value = 968.1818
$ws.Range("L34").Value = 5625.8423
$ws.Range("M34").Value = -766.1818
$ws.Range("N34").Value = -6029.8423

$ws.Range("H94").Value = 3691.5334
$ws.Range("I94").Value = 2271.75
$ws.Range("J94").Value = 5314.143
$ws.Range("K94").Value = 2271.75
$ws.Range("L94").Value = 5314.143
$ws.Range("M94").Value = -1820.75
$ws.Range("N94").Value = -6216.143

$ws.Range("H105").Value = 1206.6111
$ws.Range("I105").Value = 975.5714
$ws.Range("J105").Value = 2015.25
$ws.Range("K105").Value = 975.5714
$ws.Range("L105").Value = 2015.25
$ws.Range("M105").Value = 771.4286
$ws.Range("N105").Value = -5509.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2155.4666
$ws.Range("I3").Value = 1280.8572
$ws.Range("J3").Value = 14400
$ws.Range("K3").Value = 3842.5716
$ws.Range("L3").Value = 43200
$ws.Range("M3").Value = -3730.5716
$ws.Range("N3").Value = -43424

$ws.Range("H109").Value = 3987.5454
$ws.Range("I109").Value = 590.8
$ws.Range("J109").Value = 4986.5884
$ws.Range("K109").Value = 1772.4
$ws.Range("L109").Value = 14959.7652
$ws.Range("M109").Value = -732.3999999999999
$ws.Range("N109").Value = -17039.7652

$ws.Range("H131").Value = 739.58
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 739.58
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 2218.74
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -12298.74

$ws.Range("H141").Value = 5521.6665
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3714
$ws.Range("I80").Value = 3300
$ws.Range("J80").Value = 3873.2307
$ws.Range("K80").Value = 3300
$ws.Range("L80").Value = 3873.2307
$ws.Range("M80").Value = -2302
$ws.Range("N80").Value = -5869.2307

$ws.Range("H83").Value = 3714
$ws.Range("I83").Value = 3300
$ws.Range("J83").Value = 3873.2307
$ws.Range("K83").Value = 16500
$ws.Range("L83").Value = 19366.1535
$ws.Range("M83").Value = -11508
$ws.Range("N83").Value = -29350.1535

$ws.Range("H122").Value = 2432
$ws.Range("I122").Value = 1101.4445
$ws.Range("J122").Value = 4142.7144
$ws.Range("K122").Value = 3304.3335
$ws.Range("L122").Value = 12428.1432
$ws.Range("M122").Value = -854.3335000000002
$ws.Range("N122").Value = -17328.1432

$ws.Range("H132").Value = 41845.215
$ws.Range("I132").Value = 8745.714
$ws.Range("J132").Value = 74944.71000000001
$ws.Range("K132").Value = 26237.142
$ws.Range("L132").Value = 224834.13
$ws.Range("M132").Value = -23707.142
$ws.Range("N132").Value = -229894.13

$ws.Range("H139").Value = 25181.924
$ws.Range("J139").Value = 25181.924
$ws.Range("L139").Value = 25181.924
$ws.Range("N139").Value = -35461.924

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 3633.3333
$ws.Range("J24").Value = 3633.3333
$ws.Range("L24").Value = 3633.3333
$ws.Range("N24").Value = -4319.3333

$ws.Range("H127").Value = 38749.938
$ws.Range("J127").Value = 38749.938
$ws.Range("L127").Value = 38749.938
$ws.Range("N127").Value = -48669.938

$ws.Range("H132").Value = 2076.3635
$ws.Range("I132").Value = 1293.875
$ws.Range("J132").Value = 4163
$ws.Range("K132").Value = 3881.625
$ws.Range("L132").Value = 12489
$ws.Range("M132").Value = -1351.625
$ws.Range("N132").Value = -17549

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 71429390
$ws.Range("I107").Value = 100000340
$ws.Range("J107").Value = 2025
$ws.Range("K107").Value = 300001020
$ws.Range("L107").Value = 6075
$ws.Range("M107").Value = -299999100
$ws.Range("N107").Value = -9915

$ws.Range("H109").Value = 26980
$ws.Range("J109").Value = 26980
$ws.Range("L109").Value = 26980
$ws.Range("N109").Value = -29754
